$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "doctyp_code"
$ws.Range("C1").Value = "doccat_code"
$ws.Range("D1").Value = "is_active"

# Remove the second row (data row) entirely
$ws.Rows.Item(2).Delete()
